$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "top_movies": drop the "top_genome_tags" column (old column E), which
# shifts "short_title" (old column F) left into column E. Then refresh the
# header labels and replace the data rows with the new top-10 list.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("top_movies")

# Remove the old column E ("top_genome_tags"); everything to the right
# (short_title) shifts left automatically.
$ws1.Range("E1").EntireColumn.Delete()

# Refresh header row text (C/D swap meaning, E is now short_title).
$ws1.Range("A1").Value = "movieid"
$ws1.Range("B1").Value = "title"
$ws1.Range("C1").Value = "avg_rating"
$ws1.Range("D1").Value = "rating_count"
$ws1.Range("E1").Value = "short_title"

$top1Data = @(
    @(90464,  "Frozen North, The (2006)", 4.8, 5, "Frozen North, The (2006)"),
    @(185669, "CM Punk: Best in the World (2012)", 4.7, 5, "CM Punk: Best in the World (2012)"),
    @(150228, "Inner Worlds, Outer Worlds (2012)", 4.6, 5, "Inner Worlds, Outer Worlds (2012)"),
    @(176113, "Can't Buy My Love (2017)", 4.6, 5, "Can't Buy My Love (2017)"),
    @(171705, "Den radio (2001)", 4.58, 13, "Den radio (2001)"),
    @(104119, "Forsyte Saga, The (1967)", 4.5, 6, "Forsyte Saga, The (1967)"),
    @(173309, "Legend of the Galactic Heroes: Overture to a New War (1993)", 4.5, 6, "Legend of the Galactic Heroes: Overture to a New War (1993)"),
    @(139096, "Unmatched (2010)", 4.5, 5, "Unmatched (2010)"),
    @(91007,  "I Want to Be a Soldier (2011)", 4.5, 5, "I Want to Be a Soldier (2011)"),
    @(176569, "Liberation Day (2016)", 4.5, 5, "Liberation Day (2016)")
)

$r = 2
foreach ($row in $top1Data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Dimension / AutoFilter should now cover A1:E11 (handled automatically by
# Excel once the column is removed and data populated, but set explicitly
# to be safe).
$ws1.AutoFilterMode = $false
$ws1.Range("A1:E11").AutoFilter()

# ---------------------------------------------------------------------------
# Sheet "ratings_trend": append a new trailing row (year 2025).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("ratings_trend")

$ws3.Cells.Item(26, 1).Value = 2025
$ws3.Cells.Item(26, 2).Value = 3
$ws3.Cells.Item(26, 3).Value = 2

$ws3.AutoFilterMode = $false
$ws3.Range("A1:C26").AutoFilter()

# Extend the three per-column conditional formats from row 25 to row 26.
$rtCols = @("A", "B", "C")
foreach ($col in $rtCols) {
    $oldRange = $ws3.Range($col + "2:" + $col + "25")
    $newRange = $ws3.Range($col + "2:" + $col + "26")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($newRange)
        if ($fc.Type -eq 2) {
            $f1 = $fc.Formula1
            $f1 = $f1.Replace("`$25", "`$26")
            $fc.Formula1 = $f1
        }
    }
}

# ---------------------------------------------------------------------------
# Defined names: the hidden _FilterDatabase names must track the resized
# autofilter ranges on top_movies and ratings_trend (genres_stats is
# untouched).
# ---------------------------------------------------------------------------
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "top_movies!_FilterDatabase") {
        $n.RefersTo = "='top_movies'!`$A`$1:`$E`$11"
    } elseif ($n.Name -eq "ratings_trend!_FilterDatabase") {
        $n.RefersTo = "='ratings_trend'!`$A`$1:`$C`$26"
    }
}
